$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.039.24'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '2.039.73'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.662'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.382'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0778'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.108'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.90'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.897'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.34%  '
$ws.Range('D14').Value = '2.344.55'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').Value = '2.045.33'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.85'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +14.02%  '
$ws.Range('D18').Value = '37.113.42'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +4.77%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.08%  '
$ws.Range('E27').Value = '  -7.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.52%  '
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  +4.74%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0872'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('E39').Value = '  +10.20%  '
$ws.Range('E40').Value = '  +8.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0992'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.69%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '97.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('D47').Value = '1.282.07'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '2.227.52'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.29%  '
